$wb = $excel.ActiveWorkbook

$wsALC = $wb.Worksheets.Item("ALC")
$wsARM = $wb.Worksheets.Item("ARM")
$wsBSM = $wb.Worksheets.Item("BSM")
$wsCRP = $wb.Worksheets.Item("CRP")
$wsCUL = $wb.Worksheets.Item("CUL")
$wsGSM = $wb.Worksheets.Item("GSM")
$wsLTW = $wb.Worksheets.Item("LTW")
$wsWVR = $wb.Worksheets.Item("WVR")

# --- ALC ---
$wsALC.Range("H42").Value = 322.92856
$wsALC.Range("I42").Value = 468.42856
$wsALC.Range("J42").Value = 177.42857
$wsALC.Range("K42").Value = 1405.28568
$wsALC.Range("L42").Value = 532.28571
$wsALC.Range("M42").Value = -1175.28568
$wsALC.Range("N42").Value = -992.28571
$wsALC.Range("H74").Value = 6599.7334
$wsALC.Range("I74").Value = 5800.7144
$wsALC.Range("K74").Value = 5800.7144
$wsALC.Range("M74").Value = -4864.7144
$wsALC.Range("H77").Value = 6599.7334
$wsALC.Range("I77").Value = 5800.7144
$wsALC.Range("K77").Value = 29003.572
$wsALC.Range("M77").Value = -24323.572
$wsALC.Range("H93").Value = 31277.777
$wsALC.Range("J93").Value = 31277.777
$wsALC.Range("L93").Value = 31277.777
$wsALC.Range("N93").Value = -36269.777
$wsALC.Range("H116").Value = 440567.44
$wsALC.Range("I116").Value = 1001485.9
$wsALC.Range("J116").Value = 9091.691999999999
$wsALC.Range("K116").Value = 1001485.9
$wsALC.Range("L116").Value = 9091.691999999999
$wsALC.Range("M116").Value = -998043.9
$wsALC.Range("N116").Value = -15975.692
$wsALC.Range("H137").Value = 975258.4
$wsALC.Range("I137").Value = 2384975.2
$wsALC.Range("J137").Value = 3039.862
$wsALC.Range("K137").Value = 7154925.600000001
$wsALC.Range("L137").Value = 9119.585999999999
$wsALC.Range("M137").Value = -7152375.600000001
$wsALC.Range("N137").Value = -14219.586
$wsALC.Range("H138").Value = 2841.25
$wsALC.Range("I138").Value = 1045.125
$wsALC.Range("J138").Value = 3439.9583
$wsALC.Range("K138").Value = 3135.375
$wsALC.Range("L138").Value = 10319.8749
$wsALC.Range("M138").Value = 2004.625
$wsALC.Range("N138").Value = -20599.8749
$wsALC.Range("H141").Value = 32427.092
$wsALC.Range("I141").Value = 43398.918
$wsALC.Range("K141").Value = 130196.754
$wsALC.Range("M141").Value = -125016.754

# --- ARM ---
$wsARM.Range("H32").Value = 4642.9277
$wsARM.Range("I32").Value = 4708.8887
$wsARM.Range("J32").Value = 4405.467
$wsARM.Range("K32").Value = 4708.8887
$wsARM.Range("L32").Value = 4405.467
$wsARM.Range("M32").Value = -4421.8887
$wsARM.Range("N32").Value = -4979.467
$wsARM.Range("H97").Value = 876.125
$wsARM.Range("I97").Value = 730
$wsARM.Range("K97").Value = 730
$wsARM.Range("M97").Value = -234
$wsARM.Range("H122").Value = 3181.6155
$wsARM.Range("I122").Value = 1677.875
$wsARM.Range("J122").Value = 5587.6
$wsARM.Range("K122").Value = 5033.625
$wsARM.Range("L122").Value = 16762.8
$wsARM.Range("M122").Value = -2583.625
$wsARM.Range("N122").Value = -21662.8
$wsARM.Range("H132").Value = 3324.75
$wsARM.Range("I132").Value = 2451.8823
$wsARM.Range("K132").Value = 7355.646900000001
$wsARM.Range("M132").Value = -4825.646900000001
$wsARM.Range("H137").Value = 41770
$wsARM.Range("J137").Value = 41770
$wsARM.Range("L137").Value = 41770
$wsARM.Range("N137").Value = -51970

# --- BSM ---
$wsBSM.Range("H137").Value = 47526.668
$wsBSM.Range("J137").Value = 47526.668
$wsBSM.Range("L137").Value = 47526.668
$wsBSM.Range("N137").Value = -57726.668

# --- CRP ---
$wsCRP.Range("H50").Value = 26500.7
$wsCRP.Range("J50").Value = 26500.7
$wsCRP.Range("L50").Value = 26500.7
$wsCRP.Range("N50").Value = -27750.7
$wsCRP.Range("H51").Value = 21066.062
$wsCRP.Range("I51").Value = 8000
$wsCRP.Range("J51").Value = 21937.133
$wsCRP.Range("K51").Value = 8000
$wsCRP.Range("L51").Value = 21937.133
$wsCRP.Range("M51").Value = -7264
$wsCRP.Range("N51").Value = -23409.133
$wsCRP.Range("H60").Value = 23503.824
$wsCRP.Range("I60").Value = 20093
$wsCRP.Range("J60").Value = 23717
$wsCRP.Range("K60").Value = 20093
$wsCRP.Range("L60").Value = 23717
$wsCRP.Range("M60").Value = -19582
$wsCRP.Range("N60").Value = -24739
$wsCRP.Range("H61").Value = 21066.062
$wsCRP.Range("I61").Value = 8000
$wsCRP.Range("J61").Value = 21937.133
$wsCRP.Range("K61").Value = 8000
$wsCRP.Range("L61").Value = 21937.133
$wsCRP.Range("M61").Value = -7652
$wsCRP.Range("N61").Value = -22633.133
$wsCRP.Range("H107").Value = 643.4878
$wsCRP.Range("I107").Value = 338.55554
$wsCRP.Range("J107").Value = 1231.5714
$wsCRP.Range("K107").Value = 338.55554
$wsCRP.Range("L107").Value = 1231.5714
$wsCRP.Range("M107").Value = 1581.44446
$wsCRP.Range("N107").Value = -5071.5714

# --- CUL ---
$wsCUL.Range("H68").Value = 1291.4412
$wsCUL.Range("I68").Value = 1215.7894
$wsCUL.Range("J68").Value = 1320.7755
$wsCUL.Range("K68").Value = 3647.3682
$wsCUL.Range("L68").Value = 3962.3265
$wsCUL.Range("M68").Value = -2836.3682
$wsCUL.Range("N68").Value = -5584.3265
$wsCUL.Range("H71").Value = 1291.4412
$wsCUL.Range("I71").Value = 1215.7894
$wsCUL.Range("J71").Value = 1320.7755
$wsCUL.Range("K71").Value = 10942.1046
$wsCUL.Range("L71").Value = 11886.9795
$wsCUL.Range("M71").Value = -6886.104599999999
$wsCUL.Range("N71").Value = -19998.9795
$wsCUL.Range("H98").Value = 749.75
$wsCUL.Range("I98").Value = 633.3333
$wsCUL.Range("J98").Value = 819.6
$wsCUL.Range("K98").Value = 1899.9999
$wsCUL.Range("L98").Value = 2458.8
$wsCUL.Range("M98").Value = -401.9999
$wsCUL.Range("N98").Value = -5454.8
$wsCUL.Range("H121").Value = 1817.1692
$wsCUL.Range("I121").Value = 737.375
$wsCUL.Range("J121").Value = 1968.7192
$wsCUL.Range("K121").Value = 2212.125
$wsCUL.Range("L121").Value = 5906.1576
$wsCUL.Range("M121").Value = -902.125
$wsCUL.Range("N121").Value = -8526.1576
$wsCUL.Range("H131").Value = 767.1458
$wsCUL.Range("J131").Value = 810.5476
$wsCUL.Range("L131").Value = 2431.6428
$wsCUL.Range("N131").Value = -12511.6428
$wsCUL.Range("H132").Value = 3016.8667
$wsCUL.Range("I132").Value = 901.5
$wsCUL.Range("J132").Value = 3342.3076
$wsCUL.Range("K132").Value = 8113.5
$wsCUL.Range("L132").Value = 30080.7684
$wsCUL.Range("M132").Value = -5583.5
$wsCUL.Range("N132").Value = -35140.7684

# --- GSM ---
$wsGSM.Range("H70").Value = 6233.451
$wsGSM.Range("I70").Value = 5682.975
$wsGSM.Range("J70").Value = 8235.182000000001
$wsGSM.Range("K70").Value = 5682.975
$wsGSM.Range("L70").Value = 8235.182000000001
$wsGSM.Range("M70").Value = -5412.975
$wsGSM.Range("N70").Value = -8775.182000000001
$wsGSM.Range("H73").Value = 6233.451
$wsGSM.Range("I73").Value = 5682.975
$wsGSM.Range("J73").Value = 8235.182000000001
$wsGSM.Range("K73").Value = 5682.975
$wsGSM.Range("L73").Value = 8235.182000000001
$wsGSM.Range("M73").Value = -4746.975
$wsGSM.Range("N73").Value = -10107.182
$wsGSM.Range("H97").Value = 704.6667
$wsGSM.Range("I97").Value = 552.1429000000001
$wsGSM.Range("J97").Value = 918.2
$wsGSM.Range("K97").Value = 552.1429000000001
$wsGSM.Range("L97").Value = 918.2
$wsGSM.Range("M97").Value = -56.14290000000005
$wsGSM.Range("N97").Value = -1910.2
$wsGSM.Range("H113").Value = 1440
$wsGSM.Range("I113").Value = 1233.3334
$wsGSM.Range("J113").Value = 1750
$wsGSM.Range("K113").Value = 1233.3334
$wsGSM.Range("L113").Value = 1750
$wsGSM.Range("M113").Value = 936.6666
$wsGSM.Range("N113").Value = -6090
$wsGSM.Range("H137").Value = 39120
$wsGSM.Range("J137").Value = 48680
$wsGSM.Range("L137").Value = 48680
$wsGSM.Range("N137").Value = -58880

# --- LTW ---
$wsLTW.Range("H2").Value = 0
$wsLTW.Range("J2").Value = 0
$wsLTW.Range("L2").Value = 0
$wsLTW.Range("N2").ClearContents()
$wsLTW.Range("H21").Value = 50453.25
$wsLTW.Range("J21").Value = 50453.25
$wsLTW.Range("L21").Value = 50453.25
$wsLTW.Range("N21").Value = -50801.25
$wsLTW.Range("H100").Value = 4015.3845
$wsLTW.Range("I100").Value = 1957.1428
$wsLTW.Range("J100").Value = 6416.6665
$wsLTW.Range("K100").Value = 1957.1428
$wsLTW.Range("L100").Value = 6416.6665
$wsLTW.Range("M100").Value = -1416.1428
$wsLTW.Range("N100").Value = -7498.6665
$wsLTW.Range("H136").Value = 4444.5356
$wsLTW.Range("I136").Value = 2027.8462
$wsLTW.Range("J136").Value = 6539
$wsLTW.Range("K136").Value = 6083.5386
$wsLTW.Range("L136").Value = 19617
$wsLTW.Range("M136").Value = -3533.5386
$wsLTW.Range("N136").Value = -24717

# --- WVR ---
$wsWVR.Range("H44").Value = 40041
$wsWVR.Range("J44").Value = 40041
$wsWVR.Range("L44").Value = 40041
$wsWVR.Range("N44").Value = -41149
$wsWVR.Range("H81").Value = 1800.0667
$wsWVR.Range("I81").Value = 1336.3636
$wsWVR.Range("J81").Value = 3075.25
$wsWVR.Range("K81").Value = 2672.7272
$wsWVR.Range("L81").Value = 6150.5
$wsWVR.Range("M81").Value = -1611.7272
$wsWVR.Range("N81").Value = -8272.5
$wsWVR.Range("H84").Value = 1800.0667
$wsWVR.Range("I84").Value = 1336.3636
$wsWVR.Range("J84").Value = 3075.25
$wsWVR.Range("K84").Value = 13363.636
$wsWVR.Range("L84").Value = 30752.5
$wsWVR.Range("M84").Value = -8059.635999999999
$wsWVR.Range("N84").Value = -41360.5
$wsWVR.Range("H100").Value = 496.66666
$wsWVR.Range("J100").Value = 560
$wsWVR.Range("L100").Value = 1120
$wsWVR.Range("N100").Value = -2202
$wsWVR.Range("H132").Value = 3441
$wsWVR.Range("I132").Value = 999.8570999999999
$wsWVR.Range("J132").Value = 5149.8
$wsWVR.Range("K132").Value = 2999.5713
$wsWVR.Range("L132").Value = 15449.4
$wsWVR.Range("M132").Value = -469.5712999999996
$wsWVR.Range("N132").Value = -20509.4
